$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '58.669.01'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.71%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.151.36'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.03%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '530.78'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.78%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '139.58'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.53%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.535'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +14.55%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.31'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.12%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.433'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.61%  '
$ws.Range("E11").Value = '  +2.22%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.141'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.92%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.695.21'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.70%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.83'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.06%  '
$ws.Range("E15").Value = '  +3.57%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '58.707.77'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.69%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.25'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.40%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.149.02'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.15%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.01'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.13%  '
$ws.Range("E20").Value = '  -0.70%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '371.23'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.74%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.81'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.58%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.24%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.524'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.42%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '69.79'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.81%  '
$ws.Range("B26").Value = 'Binance-PegBSC-USD'
$ws.Range("C26").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.02'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.98%  '
$ws.Range("B27").Value = 'Kaspa'
$ws.Range("C27").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.167'
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.25'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +12.15%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0861'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.53%  '
$ws.Range("E30").Value = '  -0.21%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '22.09'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.54%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.09'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.69%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.16'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.73%  '
$ws.Range("E34").Value = '  +1.13%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '158.71'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.03%  '
$ws.Range("E36").Value = '  +2.75%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.33'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.45%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '25.10'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.21%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.68'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.03%  '
$ws.Range("E40").Value = '  +2.09%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.629.81'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.41%  '
$ws.Range("E42").Value = '  +5.58%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '38.97'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.97%  '
$ws.Range("E44").Value = '  +5.88%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.708'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.62%  '
$ws.Range("E46").Value = '  +0.05%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.194.74'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.28%  '
$ws.Range("E48").Value = '  +13.34%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.19'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.13%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.980'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.79%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '20.32'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.25%  '
